# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# that are refreshed each time the handback report is regenerated.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for the first row (d77931b3...)
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 23:10:44"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first row
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 23:10:39"
$wsZhCn.Range("K2").Value = "2016-09-04 23:10:57"

# de-de sheet: Correspond Handoff Datetime (shares its value with the
# Overview sheet's "Latest HO Xliff Generate Date") and Correspond
# Handback DateTime for the first row
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 23:10:44"
$wsDeDe.Range("K2").Value = "2016-09-04 23:11:09"
